$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Radix Sort results for each list-property block.
# Values are stored as text (matching the rest of the sheet, which uses
# text-formatted numbers), so a leading apostrophe is used to prevent
# Excel from auto-converting them into numeric values.

$ws.Range("B16").Value = "'50,000"
$ws.Range("C16").Value = "'700,000"
$ws.Range("D16").Value = "'9.97"

$ws.Range("B24").Value = "'50,000"
$ws.Range("C24").Value = "'600,000"
$ws.Range("D24").Value = "'8.30"

$ws.Range("B32").Value = "'50,000"
$ws.Range("C32").Value = "'300,000"
$ws.Range("D32").Value = "'6.28"

# Reset style back to Normal so the quote-prefix formatting introduced by
# the leading apostrophe doesn't leave a lingering custom cell style.
$ws.Range("B16:D16").Style = "Normal"
$ws.Range("B24:D24").Style = "Normal"
$ws.Range("B32:D32").Style = "Normal"
